# Generate Report for handoff
# Adds two new localization entries (0ed12709-... and 45cefd5f-...) to the
# Overview / zh-cn / de-de sheets, pushes the existing ".localization-config"
# bookkeeping row down below them, and flips the two pre-existing source
# files from "Ready for handoff" to "In Translation" (since a new handoff
# round has just been generated for the newer files).

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/e2e"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/40e3cbdba1f8e49ce88d4aae09f0e4af8dfa26b8/.localization-config"
$zhcnBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0920e71b8e6468777c03d1a93dacdbeffb2766ed/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu"
$dedeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d61823019d4254cb5dc26a4b774f1312884e132/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu"

$file3 = "0ed12709-d088-4d8c-8475-0e19163a68f3"
$file4 = "45cefd5f-f1e5-46e1-9604-137c12761e97"
$hash3 = "bc16ce64bd8926fc2a9dfb5ce635faaa02561a40"
$hash4 = "ff005cdba3b917e92e3a7ff3db8ae7ae7f76afa5"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Existing two files move from "Ready for handoff" to "In Translation".
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"
$ov.Range("B3").Value = "In Translation"
$ov.Range("C3").Value = "In Translation"

# New rows 4 and 5 for the newly generated handoff files.
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"

# Row 6 is the relocated ".localization-config" bookkeeping row.
$ov.Range("B6").Value = "Not to be localized"
$ov.Range("C6").Value = "Not to be localized"

# Hyperlinks: this engine can only clear *all* hyperlinks on a sheet at once,
# so drop them all and recreate the full, final set in order.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "$mdBase/88b63550-690b-4866-9085-6bae40f52bf6.md", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "$mdBase/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "$mdBase/$file3.md", "", "", "$file3.md")
$ov.Hyperlinks.Add($ov.Range("A5"), "$mdBase/$file4.md", "", "", "$file4.md")
$ov.Hyperlinks.Add($ov.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = "In Translation"
$zh.Range("B3").Value = "In Translation"

$zh.Range("B4").Value = "Ready for handoff"
$zh.Range("D4").Value = "2016-01-25 03:06:42"
$zh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Include"

$zh.Range("A5").Value = "$file4.md"
$zh.Range("B5").Value = "Ready for handoff"
$zh.Range("D5").Value = "2016-01-25 03:06:42"
$zh.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G5").Value = "0001-01-01 00:00:00"
$zh.Range("H5").Value = "Include"

$zh.Range("A6").Value = ".localization-config"
$zh.Range("B6").Value = "Not to be localized"
$zh.Range("D6").Value = "0001-01-01 00:00:00"
$zh.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G6").Value = "0001-01-01 00:00:00"
$zh.Range("H6").Value = "Ignored"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "$mdBase/88b63550-690b-4866-9085-6bae40f52bf6.md", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "$zhcnBase/88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.zh-cn.xlf", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "$mdBase/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "$zhcnBase/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.zh-cn.xlf", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "$mdBase/$file3.md", "", "", "$file3.md")
$zh.Hyperlinks.Add($zh.Range("C4"), "$zhcnBase/$file3.$hash3.zh-cn.xlf", "", "", "$file3.$hash3.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A5"), "$mdBase/$file4.md", "", "", "$file4.md")
$zh.Hyperlinks.Add($zh.Range("C5"), "$zhcnBase/$file4.$hash4.zh-cn.xlf", "", "", "$file4.$hash4.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = "In Translation"
$de.Range("B3").Value = "In Translation"

$de.Range("B4").Value = "Ready for handoff"
$de.Range("D4").Value = "2016-01-25 03:06:53"
$de.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Include"

$de.Range("A5").Value = "$file4.md"
$de.Range("B5").Value = "Ready for handoff"
$de.Range("D5").Value = "2016-01-25 03:06:53"
$de.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G5").Value = "0001-01-01 00:00:00"
$de.Range("H5").Value = "Include"

$de.Range("A6").Value = ".localization-config"
$de.Range("B6").Value = "Not to be localized"
$de.Range("D6").Value = "0001-01-01 00:00:00"
$de.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G6").Value = "0001-01-01 00:00:00"
$de.Range("H6").Value = "Ignored"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "$mdBase/88b63550-690b-4866-9085-6bae40f52bf6.md", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.md")
$de.Hyperlinks.Add($de.Range("C2"), "$dedeBase/88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.de-de.xlf", "", "", "88b63550-690b-4866-9085-6bae40f52bf6.7187555e9660a7463e2b9e7be747327a8f8f343d.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "$mdBase/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.md")
$de.Hyperlinks.Add($de.Range("C3"), "$dedeBase/eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.de-de.xlf", "", "", "eaff493e-b3d4-4f4e-90b1-64f9b3f29d33.34f728b61c45b280c81aa4a7d18264a462f6403a.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "$mdBase/$file3.md", "", "", "$file3.md")
$de.Hyperlinks.Add($de.Range("C4"), "$dedeBase/$file3.$hash3.de-de.xlf", "", "", "$file3.$hash3.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A5"), "$mdBase/$file4.md", "", "", "$file4.md")
$de.Hyperlinks.Add($de.Range("C5"), "$dedeBase/$file4.$hash4.de-de.xlf", "", "", "$file4.$hash4.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A6"), $cfgUrl, "", "", ".localization-config")
